# Actualización automática 2025-06-09 15:45:08
# Inserts a new client row ("JUAREZ FLORES JORGE WILLIAMS") right before
# "LUI WONG ANGEL BOLIVAR" (row 34) on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, shifting the existing rows down by one, and
# updates the trailing "N de 55" -> "N de 56" summary labels on
# "VENTAS POR GRUPO".

$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO (columns A:N) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(34).Insert()
$ws1.Range("A34").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Range("B34").Value = "JUAREZ FLORES JORGE WILLIAMS"
$ws1.Range("C34:N34").Value = 0

# The summary row (previously row 57, now row 58) counts "de 55" -> "de 56"
$summaryCols1 = @("C","D","E","F","G","H","I","J","K","L","M","N")
foreach ($col in $summaryCols1) {
    $cell = $ws1.Range($col + "58")
    $cur = $cell.Value()
    $cell.Value = $cur -replace "de 55", "de 56"
}

# ---- Sheet: VENTA MENSUAL (columns A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(34).Insert()
$ws2.Range("A34").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Range("B34").Value = "JUAREZ FLORES JORGE WILLIAMS"
$ws2.Range("C34:G34").Value = 0
